# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Updates Price (col D) and Volume(1h) (col E) for existing rows, and
# reorders three rows (Dai/OKB/Maker rotate, Fetch.AI/Stacks swap) which
# changes their Coin (B), Link (C), Price (D) and Volume(1h) (E) values.
#
# Some new Price values parse as plain numbers (e.g. "580.95"); those are
# prefixed with a leading apostrophe so Excel stores them as literal text
# (matching the original column, which holds text like "67.777.65" or
# "0.999") instead of silently converting them to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.777.65"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "3.336.15"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'580.95"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").Value = "'176.23"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +2.54%  "
$ws.Range("D9").Value = "3.332.68"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("E10").Value = "  +7.24%  "
$ws.Range("E11").Value = "  +2.87%  "
$ws.Range("D12").Value = "'47.06"
$ws.Range("E12").Value = "  +5.81%  "
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("D14").Value = "'693.41"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "3.875.63"
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").Value = "'8.44"
$ws.Range("E16").Value = "  +3.45%  "
$ws.Range("D17").Value = "67.807.69"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "3.346.85"
$ws.Range("E19").Value = "  +2.59%  "
$ws.Range("D20").Value = "'17.54"
$ws.Range("E20").Value = "  +2.91%  "
$ws.Range("D21").Value = "'11.08"
$ws.Range("E21").Value = "  +4.88%  "
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").Value = "'5.46"
$ws.Range("E23").Value = "  +6.04%  "
$ws.Range("D24").Value = "'16.98"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").Value = "'100.77"
$ws.Range("E25").Value = "  +3.36%  "
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("E27").Value = "  +3.40%  "
$ws.Range("D28").Value = "'9.56"
$ws.Range("E28").Value = "  +6.51%  "
$ws.Range("D29").Value = "'33.05"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").Value = "'8.57"
$ws.Range("E30").Value = "  +4.43%  "
$ws.Range("D31").Value = "'7.08"
$ws.Range("E31").Value = "  +8.46%  "
$ws.Range("D32").Value = "'565.05"
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("E34").Value = "  +3.90%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'57.30"
$ws.Range("E35").Value = "  +3.93%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.714.74"
$ws.Range("E36").Value = "  -2.61%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "'3.34"
$ws.Range("E38").Value = "  +2.78%  "
$ws.Range("D39").Value = "'35.02"
$ws.Range("E39").Value = "  +12.64%  "
$ws.Range("E40").Value = "  +5.60%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'3.17"
$ws.Range("E41").Value = "  +7.68%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'2.63"
$ws.Range("E42").Value = "  +3.56%  "
$ws.Range("E43").Value = "  +3.53%  "
$ws.Range("E44").Value = "  +4.99%  "
$ws.Range("D45").Value = "'3.28"
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("E46").Value = "  +2.86%  "
$ws.Range("E47").Value = "  +6.36%  "
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "'1.33"
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("D51").Value = "'131.67"
$ws.Range("E51").Value = "  +2.97%  "
